# Update cryptocurrency price/volume data in the active worksheet
# (values refreshed by the scheduled GitHub Actions data pull).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into column D (Price) while preserving it as text,
# since these price strings (e.g. "29.386.62", "1.000") must not be
# reinterpreted by Excel as numbers. The original cell style is restored
# afterwards so no extra formatting is left behind.
function Set-PriceText($cellRef, $text) {
    $c = $ws.Range($cellRef)
    $origStyle = $c.Style
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.Style = $origStyle
}

Set-PriceText "D2" "29.386.62"
$ws.Range("E2").Value = "  -0.01%  "
Set-PriceText "D3" "1.848.15"
$ws.Range("E3").Value = "  +0.00%  "
Set-PriceText "D5" "240.32"
$ws.Range("E5").Value = "  -0.06%  "
Set-PriceText "D6" "0.6283"
$ws.Range("E6").Value = "  -0.19%  "
Set-PriceText "D7" "1.000"
$ws.Range("E7").Value = "  +0.00%  "
Set-PriceText "D8" "0.07639"
$ws.Range("E8").Value = "  +0.21%  "
$ws.Range("E9").Value = "  -1.11%  "
Set-PriceText "D10" "24.75"
$ws.Range("E10").Value = "  +0.92%  "
Set-PriceText "D11" "0.07739"
$ws.Range("E11").Value = "  -0.02%  "
Set-PriceText "D12" "5.033"
$ws.Range("E12").Value = "  +0.54%  "
Set-PriceText "D13" "0.6791"
$ws.Range("E13").Value = "  +0.01%  "
Set-PriceText "D14" "0.00001061"
$ws.Range("E14").Value = "  -2.41%  "
Set-PriceText "D15" "83.26"
$ws.Range("E15").Value = "  -0.43%  "
Set-PriceText "D16" "6.152"
$ws.Range("E16").Value = "  -0.03%  "
Set-PriceText "D17" "29.408.67"
$ws.Range("E17").Value = "  -0.02%  "
Set-PriceText "D18" "227.41"
$ws.Range("E18").Value = "  -0.62%  "
Set-PriceText "D19" "12.34"
$ws.Range("E19").Value = "  -0.72%  "
Set-PriceText "D20" "0.9998"
$ws.Range("E20").Value = "  -0.04%  "
Set-PriceText "D21" "7.502"
$ws.Range("E21").Value = "  +0.59%  "
Set-PriceText "D22" "1.000"
$ws.Range("E22").Value = "  -0.02%  "
Set-PriceText "D23" "158.32"
$ws.Range("E23").Value = "  +0.42%  "
Set-PriceText "D24" "0.1384"
$ws.Range("E24").Value = "  -0.28%  "
Set-PriceText "D25" "8.406"
Set-PriceText "D26" "17.69"
$ws.Range("E26").Value = "  +0.32%  "
Set-PriceText "D27" "1.377"
$ws.Range("E27").Value = "  +4.84%  "
Set-PriceText "D28" "1.458"
$ws.Range("E28").Value = "  -0.58%  "
Set-PriceText "D29" "0.05601"
$ws.Range("E29").Value = "  -0.46%  "
Set-PriceText "D30" "4.118"
$ws.Range("E30").Value = "  +0.16%  "
Set-PriceText "D31" "4.077"
$ws.Range("E31").Value = "  +1.01%  "
Set-PriceText "D32" "1.837"
$ws.Range("E32").Value = "  -0.75%  "
Set-PriceText "D33" "1.162"
$ws.Range("E33").Value = "  +0.50%  "
Set-PriceText "D34" "0.6942"
$ws.Range("E34").Value = "  -2.11%  "
Set-PriceText "D35" "2.578"
$ws.Range("E35").Value = "  -0.20%  "
$ws.Range("E36").Value = "  +0.18%  "
Set-PriceText "D37" "1.229.88"
$ws.Range("E37").Value = "  -0.27%  "
Set-PriceText "D38" "2.714"
$ws.Range("E38").Value = "  -2.21%  "
Set-PriceText "D39" "6.386"
$ws.Range("E39").Value = "  -1.11%  "
Set-PriceText "D40" "0.9049"
$ws.Range("E40").Value = "  -0.33%  "
Set-PriceText "D41" "1.000"
$ws.Range("E41").Value = "  +0.05%  "
Set-PriceText "D42" "101.55"
$ws.Range("E42").Value = "  -0.01%  "
Set-PriceText "D43" "66.07"
$ws.Range("E43").Value = "  +0.10%  "
Set-PriceText "D44" "7.176"
$ws.Range("E44").Value = "  +0.24%  "
Set-PriceText "D45" "0.4011"
$ws.Range("E45").Value = "  -0.09%  "
Set-PriceText "D46" "8.982"
$ws.Range("E46").Value = "  +0.00%  "
$ws.Range("E47").Value = "  -0.73%  "
Set-PriceText "D48" "0.1139"
$ws.Range("E48").Value = "  +1.41%  "
Set-PriceText "D49" "0.05701"
$ws.Range("E49").Value = "  -0.15%  "
Set-PriceText "D50" "0.4628"
$ws.Range("E50").Value = "  +0.10%  "
$ws.Range("B51").Value = "SynthetixNetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
Set-PriceText "D51" "2.532"
$ws.Range("E51").Value = "  +0.00%  "
